# edit.ps1 -- applies the "Professional Summary" rewrite + minor
# indentation/cleanup changes described by the target diff.

$d = $word.ActiveDocument

function Insert-RunXml($Range, $InnerBodyXml) {
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $InnerBodyXml + '</w:p></w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $Range.InsertXML($pkg)
}

# ---------------------------------------------------------------------
# 1) Professional summary paragraphs: the two paragraphs swap roles.
#    Old paragraph 4 ("I am looking for a job...") becomes the old
#    paragraph 5 text ("Experienced Scientific Programmer...").
#    Old paragraph 5 becomes a reworded / multi-run version of the
#    original paragraph-4 text.
# ---------------------------------------------------------------------

$p4 = $d.Paragraphs.Item(4)
$p4Range = $d.Range($p4.Range.Start, $p4.Range.End - 1)
Insert-RunXml $p4Range '<w:r><w:rPr/><w:t>Experienced Scientific Programmer with a demonstrated history of working in the Computer Software industry. Skilled in Continuous Improvement, Machine Learning, Optimization. Strong engineering professional with a Bachelor of Science, Master of Engineering, Doctor of Philosophy (Uncompleted) all focused in Industrial Engineering from Rochester Institute of Technology.</w:t></w:r>'

# paragraph 5 is unaffected in position/count, but re-fetch to be safe
$p5 = $d.Paragraphs.Item(5)
$p5Range = $d.Range($p5.Range.Start, $p5.Range.End - 1)

$p5Runs = ''
$p5Runs += '<w:r><w:rPr/><w:t>Search</w:t></w:r>'
$p5Runs += '<w:r><w:rPr/><w:t xml:space="preserve">ing for a job where my analytical skills can be a part of contemporary solutions. In this job, computer mathematical models </w:t></w:r>'
$p5Runs += '<w:r><w:rPr/><w:t xml:space="preserve">would be built </w:t></w:r>'
$p5Runs += '<w:r><w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="" w:cstheme="minorBidi" w:eastAsiaTheme="minorHAnsi"/><w:color w:val="auto"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t>with</w:t></w:r>'
$p5Runs += '<w:r><w:rPr/><w:t xml:space="preserve"> knowledge from </w:t></w:r>'
$p5Runs += '<w:r><w:rPr/><w:t>data,</w:t></w:r>'
$p5Runs += '<w:r><w:rPr/><w:t xml:space="preserve"> people, </w:t></w:r>'
$p5Runs += '<w:r><w:rPr/><w:t xml:space="preserve">and </w:t></w:r>'
$p5Runs += '<w:r><w:rPr/><w:t>literature.</w:t></w:r>'

Insert-RunXml $p5Range $p5Runs

# ---------------------------------------------------------------------
# 2) Two "Subtitle" paragraphs (Researcher (Student) / RIT Sep 2016 -
#    Nov 2018 header) get an explicit w:left="0" added to their
#    (already w:hanging="0") indentation.
# ---------------------------------------------------------------------

foreach ($idx in 26, 27) {
    $p = $d.Paragraphs.Item($idx)
    $p.Format.LeftIndent = 0
}

# ---------------------------------------------------------------------
# 3) Remove the three empty trailing "Subtitle" paragraphs that follow
#    the "Modeled the distinguishing characteristics of bladder cancer
#    patients..." bullet (keep the following "Data Analyst (Intern)"
#    paragraph intact).
# ---------------------------------------------------------------------

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match 'bladder cancer') {
        $target = $p
    }
}

if ($target -ne $null) {
    $deleteStart = $target.Range.End
    $p1 = $target.Next()
    $p2 = $p1.Next()
    $p3 = $p2.Next()
    $deleteEnd = $p3.Range.End
    $r = $d.Range($deleteStart, $deleteEnd)
    $r.Delete()
}

Write-Output "edit.ps1 complete"
